$p = $ppt.ActivePresentation
Write-Output ("Before: " + $p.PageSetup.SlideWidth + " x " + $p.PageSetup.SlideHeight)
$p.PageSetup.SlideWidth = 11520488 / 12700.0
$p.PageSetup.SlideHeight = 3240088 / 12700.0
Write-Output ("After: " + $p.PageSetup.SlideWidth + " x " + $p.PageSetup.SlideHeight)
